# Update "想去人数" (want-to-go count) figures to the newly scraped values.
# Sheet order in the workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsPerformance = $wb.Worksheets.Item(2)  # 演出
$wsAllTypes = $wb.Worksheets.Item(4)     # 全部类型

# 展览 sheet
$wsExhibition.Range("F5").Value = 13700
$wsExhibition.Range("F15").Value = 13697
$wsExhibition.Range("F37").Value = 3710

# 演出 sheet
$wsPerformance.Range("F2").Value = 44

# 全部类型 sheet (combined view of 展览 + 演出)
$wsAllTypes.Range("F5").Value = 13700
$wsAllTypes.Range("F15").Value = 13697
$wsAllTypes.Range("F31").Value = 44
$wsAllTypes.Range("F39").Value = 3712
